$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "company_code"
$ws.Range("B1").Value = "branch_code"
$ws.Range("C1").Value = "category_code"
$ws.Range("D1").Value = "code"
$ws.Range("E1").Value = "name"
$ws.Range("F1").Value = "description"

# Row 2
$ws.Range("A2").Value = "C001"
$ws.Range("B2").Value = "B001"
$ws.Range("C2").Value = "CAT001"
$ws.Range("D2").Value = "SUBC001"
$ws.Range("E2").Value = "sub category pertama"
$ws.Range("F2").Value = "description  sub category pertama"

# Row 3
$ws.Range("A3").Value = "C002"
$ws.Range("B3").Value = "B002"
$ws.Range("C3").Value = "CAT002"
$ws.Range("D3").Value = "SUBC002"
$ws.Range("E3").Value = "sub categoroy kedua"
$ws.Range("F3").Value = "description sub category kedua"
